$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(219).Delete()
